$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-19 06:51:38"
$wsZhCn.Range("H2").Value = "2016-03-19 06:51:59"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-19 06:51:41"
$wsDeDe.Range("H2").Value = "2016-03-19 06:52:04"
